$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45913
$ws.Range("B2").Value = 102.3
$ws.Range("C2").Value = 99.54000000000001
$ws.Range("D2").Value = 97.95999999999999
$ws.Range("E2").Value = 97.84
$ws.Range("F2").Value = 96.68000000000001
$ws.Range("G2").Value = 92
$ws.Range("H2").Value = 94.98
$ws.Range("I2").Value = 100.02
$ws.Range("J2").Value = 99.78
$ws.Range("K2").Value = 69.04000000000001
$ws.Range("L2").Value = 9.619999999999999
$ws.Range("M2").Value = 2.37
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 0
$ws.Range("S2").Value = 0.65
$ws.Range("T2").Value = 14
$ws.Range("U2").Value = 82.86
$ws.Range("V2").Value = 100.13
$ws.Range("W2").Value = 106.11
$ws.Range("X2").Value = 99.3
$ws.Range("Y2").Value = 95.09999999999999
$ws.Range("Z2").Value = 60.85
$ws.Range("AB2").Value = 100.16
$ws.Range("AD2").Value = 103.12
$ws.Range("AE2").Value = "0h-2h"
$ws.Range("AF2").Value = 100.92
